$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N/O/P -> O/P/Q)
[void]$ws.Columns("N").Insert()

# New column N takes the same display width as column M (~11 chars stored width)
$ws.Columns("N").ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet / tab, with the new selection
$ws.Activate()
[void]$ws.Range("R6").Select()
